$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Description")

# Remove the extraneous "If there are increases..." note row. Deleting the
# entire row shifts the rows below it up by one (matching rows 9-13 -> 8-12)
# and removes the now-unused shared string from the workbook's string table.
$ws.Rows.Item(8).Delete()

# Reflect the resulting active cell/selection (as left after the row delete).
$ws.Range("A8").Select()
